# The hotel room-price table has a footnote paragraph reading "* Chinese
# only" attached to the "Double Room" row. Remove that note text (the
# paragraph mark / empty paragraph itself stays in place).
$d = $word.ActiveDocument

$d.Content.Find.Execute("* Chinese only", $true, $false, $false, $false, $false,
                         $true, 1, $false, "", 2)
